$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the flow-rate values for the third data block (35 -> 50 ml/min)
$ws.Range("A18:A25").Value = 50

# Remove the intermediate unit-conversion columns (old B: m3/s, old C: s/m3).
# This shifts old column D (Particle Size) into B and old column E (Potential) into C.
$ws.Range("B:C").Delete()

# Re-anchor the floating (empty) text boxes that used to sit just inside the
# old column C onto the new column B (same relative offset within the column).
for ($i = 1; $i -le $ws.Shapes.Count; $i++) {
    $shp = $ws.Shapes.Item($i)
    $shp.Left = $ws.Range("B1").Left
}

$ws.Range("E8").Select()
